# Update the "Förändrad" (changed) date column C for rows 2-11
# from 2023-10-22 (serial 45221) to 2023-10-25 (serial 45224).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
